# plantilla_de_carga_NI.xlsx - "carga masiva de consumo, incluyendo solicitante y area"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns: SOLICITANTE (F1) and AREA (G1) ---------------------------
# (existing style s="1" on F1/G1 is preserved automatically; only the value
# is being set, turning the blank numeric cells into shared-string cells)
$ws.Range("F1").Value = "SOLICITANTE"
$ws.Range("G1").Value = "AREA"

# --- Column widths -----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 23.666666666666668
$ws.Columns.Item(2).ColumnWidth = 22.666666666666668
$ws.Columns.Item(3).ColumnWidth = 28.333333333333332
$ws.Columns.Item(4).ColumnWidth = 26.833333333333332
$ws.Columns.Item(5).ColumnWidth = 21.5
$ws.Columns.Item(6).ColumnWidth = 13.666666666666666

# --- Header cell comments (instructions for each field) -----------------
$ws.Range("A1").AddComment("Obligatorio`n")
$ws.Range("B1").AddComment("En caso de no tener serie dejar en blanco, evitar espacios en blanco`n")
$ws.Range("C1").AddComment("Obligatorio")
$ws.Range("D1").AddComment("Campo opcional, tener en cuenta las maquinas registradas en el contrato escribiendo el NOMBRE, de caso contrario validara con un error.`n")
$ws.Range("E1").AddComment("Opcional. Para no tener inconvenientes en el reporte por vales escribir de la siguiente manera`n(VS N° 066174) respetando los espacios")
$ws.Range("F1").AddComment("CAMPO OPCIONAL,`nINDICAR EL NUMERO DE SOLICITANTE CON EL QUE SE REGISTRO EN LA APLICACIÓN,DE LO CONTRARIO NO REGISTRARA EL CONSUMO")
$ws.Range("G1").AddComment("Campo opcional,`nIndicar el numero asignado al area en el aplicativo para evitar errores`n")

# --- Page setup (printing orientation) -----------------------------------
$ws.PageSetup.Orientation = 1

# --- Selection left at E8 by the author before saving --------------------
$ws.Range("E8").Select()
